$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.379.44'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.19%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.573.96'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.68%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -1.36%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.001'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.50%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '290.52'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.42%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3762'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.72%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '50.20'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.74%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3420'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.41%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.170'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.11%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07695'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.21%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.32%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.40'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.05%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.010'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.45%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.927'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.29%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001146'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.07%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.571.87'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.20%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.39'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.57%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06738'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.30%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.58%  '

# Row 21
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.81'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.53%  '

# Row 22
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.256'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.43%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.5310'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -8.11%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.02'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.77%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '22.387.59'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.21%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.395'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.10%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.792'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.34%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.32'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.98%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '144.83'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.24%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.075'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.60%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '126.71'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.90%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.747.95'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.87%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.037'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +13.77%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.287'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.34%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.044'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.61%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.15'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.11%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.08550'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.50%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02559'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.43%  '

# Row 39
$ws.Range("E39").Value = '  +2.68%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.545'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.49%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.06505'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.51%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.305'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.56%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.70'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.20%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6447'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.67%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.10'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.98%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.001'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.52%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6052'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.69%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.775'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.03%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.301'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +11.52%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.103'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.31%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '125.34'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.32%  '
